$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cricosCode/department column (C2) held "ACE AVIATION" and is being
# corrected to just "AVIATION".
$ws.Range("C2").Value = "AVIATION"
